$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1, matching style of existing header row (bold/bordered/centered)
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "intervention_type"

# intervention_type values for each clinical trial row (K2:K121)
$interventionTypes = @{
  2 = "PROCEDURE"
  3 = "DRUG"
  4 = "BIOLOGICAL"
  5 = "PROCEDURE"
  6 = "OTHER"
  7 = "OTHER"
  8 = "DRUG"
  9 = "DRUG"
  10 = "OTHER"
  11 = "DRUG"
  12 = "OTHER"
  13 = "DRUG"
  14 = "RADIATION"
  15 = "DRUG"
  16 = "OTHER"
  17 = "OTHER"
  18 = "BIOLOGICAL"
  19 = "DRUG"
  20 = "BIOLOGICAL"
  21 = "PROCEDURE"
  23 = "BIOLOGICAL"
  24 = "OTHER"
  25 = "BIOLOGICAL"
  26 = "PROCEDURE"
  27 = "PROCEDURE"
  28 = "PROCEDURE"
  29 = "DIETARY_SUPPLEMENT"
  30 = "DRUG"
  31 = "DRUG"
  32 = "OTHER"
  33 = "DRUG"
  34 = "DRUG"
  35 = "PROCEDURE"
  36 = "DRUG"
  37 = "OTHER"
  39 = "OTHER"
  40 = "OTHER"
  41 = "OTHER"
  42 = "OTHER"
  43 = "DRUG"
  44 = "BIOLOGICAL"
  45 = "DEVICE"
  46 = "DRUG"
  47 = "OTHER"
  48 = "OTHER"
  49 = "DRUG"
  50 = "OTHER"
  51 = "OTHER"
  52 = "OTHER"
  53 = "DIAGNOSTIC_TEST"
  54 = "BIOLOGICAL"
  55 = "DIETARY_SUPPLEMENT"
  56 = "PROCEDURE"
  57 = "BIOLOGICAL"
  59 = "DRUG"
  60 = "DRUG"
  61 = "BIOLOGICAL"
  62 = "OTHER"
  63 = "BIOLOGICAL"
  64 = "OTHER"
  65 = "OTHER"
  66 = "DRUG"
  67 = "OTHER"
  68 = "BIOLOGICAL"
  70 = "OTHER"
  71 = "PROCEDURE"
  72 = "OTHER"
  73 = "BIOLOGICAL"
  74 = "OTHER"
  75 = "OTHER"
  76 = "OTHER"
  77 = "DEVICE"
  78 = "DRUG"
  79 = "OTHER"
  80 = "BIOLOGICAL"
  81 = "OTHER"
  82 = "RADIATION"
  83 = "BIOLOGICAL"
  84 = "BIOLOGICAL"
  85 = "OTHER"
  86 = "DEVICE"
  87 = "PROCEDURE"
  88 = "OTHER"
  89 = "BIOLOGICAL"
  90 = "OTHER"
  91 = "BIOLOGICAL"
  92 = "PROCEDURE"
  93 = "OTHER"
  94 = "BIOLOGICAL"
  95 = "DRUG"
  96 = "DRUG"
  97 = "OTHER"
  98 = "PROCEDURE"
  99 = "OTHER"
  100 = "DRUG"
  101 = "BEHAVIORAL"
  102 = "DEVICE"
  103 = "OTHER"
  104 = "OTHER"
  105 = "OTHER"
  106 = "DRUG"
  107 = "OTHER"
  108 = "BIOLOGICAL"
  109 = "DIETARY_SUPPLEMENT"
  110 = "DEVICE"
  111 = "BIOLOGICAL"
  112 = "BIOLOGICAL"
  113 = "OTHER"
  114 = "BIOLOGICAL"
  115 = "OTHER"
  116 = "BIOLOGICAL"
  117 = "DEVICE"
  118 = "DRUG"
  119 = "OTHER"
  120 = "BIOLOGICAL"
}

foreach ($row in $interventionTypes.Keys) {
  $ws.Cells.Item($row, 11).Value = $interventionTypes[$row]
}

# Rows without a matched NCT record have no intervention_type data;
# still materialize the K cell (present but blank) to match the source rows.
$emptyRows = @(22, 38, 58, 69, 121)
foreach ($row in $emptyRows) {
  $ws.Cells.Item($row, 11).NumberFormat = "@"
}